# Repull data, push all data, mean calculation
# Updates specific values in column F ("dSF") to reflect repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -1
    8  = -1
    13 = 5
    14 = 2
    15 = 5
    19 = -1
    20 = 0
    21 = 2
    23 = 1
    24 = -1
    25 = 2
    30 = -1
    32 = -1
    38 = 1
    42 = -1
    44 = -5
    45 = -3
    54 = 6
    55 = -2
    58 = -2
    61 = 1
    65 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
